# "Att até Jogo 46"
# Enter the result of Jogo 46 (Juventus x Manchester City) and Jogo 45
# (Wydad Casablanca x Al Ain) on the "Fase de Grupos" sheet. This completes
# the Group G round-robin; every standings/ranking formula in that block
# (and the downstream "Finais" lookups for 1st/2nd of Group G) recalculates
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fase de Grupos")

# Jogo 46: Juventus 2 x 5 Manchester City
$ws.Range("F43").Value = 2
$ws.Range("H43").Value = 5

# Jogo 45: Wydad Casablanca 1 x 2 Al Ain
$ws.Range("F44").Value = 1
$ws.Range("H44").Value = 2

# Move the active selection to where the author left off (F49), matching
# the saved sheetView selection in the workbook.
$ws.Range("F49").Select()
